# Update countries & provincias Spain
# Applies the daily COVID data refresh to the "Pais" sheet:
#  - refreshes case/recovered/death counters for a set of countries
#  - two pairs of adjacent countries swap rank (their row positions trade
#    values because the update changed their relative ordering)
#  - bumps the "last updated" timestamp string

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Estados Unidos (row 4)
Set-Row 4 8049854 12065 5197125 2632448 0 270 220281

# Reino Unido (row 15)
Set-Row 15 634920 17234 0 0 0 143 43018

# Italia (row 20)
Set-Row 20 365467 5901 242028 87193 0 41 36246

# Turquia (row 24)
Set-Row 24 338779 1632 296972 32850 0 62 8957

# Israel (row 27)
Set-Row 27 295625 1594 242027 51566 0 11 2032

# Chequia (row 38)
Set-Row 38 125328 3907 58417 65814 0 46 1097

# Singapur (row 62) - only D & E change
$ws.Cells.Item(62, 4).Value = 57740
$ws.Cells.Item(62, 5).Value = 116

# Rows 83/84 swap order: Jordania overtakes Australia
$ws.Cells.Item(83, 1).Value = "Jordania"
Set-Row 83 28127 2054 6361 21541 0 18 225

$ws.Cells.Item(84, 1).Value = "Australia"
Set-Row 84 27316 30 25037 1380 0 1 899

# Guayana Francesa (row 108) - B,C,D,E change
$ws.Cells.Item(108, 2).Value = 10192
$ws.Cells.Item(108, 3).Value = 12
$ws.Cells.Item(108, 4).Value = 9881
$ws.Cells.Item(108, 5).Value = 242

# Malaui (row 123)
Set-Row 123 5827 3 4688 958 0 1 181

# Rows 133/134 swap order: Sri Lanka overtakes Ruanda
$ws.Cells.Item(133, 1).Value = "Sri Lanka"
Set-Row 133 5038 194 3328 1697 0 0 13

$ws.Cells.Item(134, 1).Value = "Ruanda"
Set-Row 134 4905 0 3877 996 0 0 32

# Rows 159/160 swap order: Republica de Chipre overtakes Yemen
$ws.Cells.Item(159, 1).Value = "Republica de Chipre"
Set-Row 159 2109 62 1444 640 0 0 25

$ws.Cells.Item(160, 1).Value = "Yemen"
Set-Row 160 2052 0 1329 127 0 0 596

# Liberia (row 165) - only D & E change
$ws.Cells.Item(165, 4).Value = 1251
$ws.Cells.Item(165, 5).Value = 38

# Curazao (row 174) - B,C,D,E change
$ws.Cells.Item(174, 2).Value = 619
$ws.Cells.Item(174, 3).Value = 34
$ws.Cells.Item(174, 4).Value = 343
$ws.Cells.Item(174, 5).Value = 275

# Bump the "last updated" timestamp (row 1, column A)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 13 de Octubre de 2020 a las 18:31"
